$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("N33").Value = -983
$ws.Range("I33").Value = 219.21428
$ws.Range("K33").Value = 219.21428
$ws.Range("H33").Value = 257.4375
$ws.Range("M33").Value = 9.785719999999998
$ws.Range("L33").Value = 525
$ws.Range("J33").Value = 525
$ws.Range("I40").Value = 980.36365
$ws.Range("N40").Value = -1431.6
$ws.Range("M40").Value = -805.36365
$ws.Range("H40").Value = 1012
$ws.Range("L40").Value = 1081.6
$ws.Range("J40").Value = 1081.6
$ws.Range("K40").Value = 980.36365
$ws.Range("K74").Value = 4875
$ws.Range("H74").Value = 4750
$ws.Range("M74").Value = -3939
$ws.Range("I74").Value = 4875
$ws.Range("I77").Value = 4875
$ws.Range("H77").Value = 4750
$ws.Range("K77").Value = 24375
$ws.Range("M77").Value = -19695
$ws.Range("I116").Value = 2227.2727
$ws.Range("L116").Value = 2900
$ws.Range("K116").Value = 2227.2727
$ws.Range("J116").Value = 2900
$ws.Range("M116").Value = 1214.7273
$ws.Range("H116").Value = 2547.6191
$ws.Range("N116").Value = -9784
$ws.Range("H119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L121").Value = 4350
$ws.Range("J121").Value = 1450
$ws.Range("H121").Value = 1450
$ws.Range("K121").Value = 0
$ws.Range("N121").Value = -7844
$ws.Range("I121").Value = 0
$ws.Range("H132").Value = 289419.72
$ws.Range("L132").Value = 4080
$ws.Range("K132").Value = 1012288.98
$ws.Range("M132").Value = -1009758.98
$ws.Range("I132").Value = 337429.66
$ws.Range("N132").Value = -9140
$ws.Range("J132").Value = 1360
$ws.Range("J137").Value = 6501
$ws.Range("M137").Value = -1357188.75
$ws.Range("I137").Value = 453246.25
$ws.Range("H137").Value = 404949.47
$ws.Range("L137").Value = 19503
$ws.Range("N137").Value = -24603
$ws.Range("K137").Value = 1359738.75
$ws.Range("N119").ClearContents()
$ws.Range("M121").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("J68").Value = 63600
$ws.Range("H68").Value = 63600
$ws.Range("N68").Value = -65222
$ws.Range("L68").Value = 63600
$ws.Range("H71").Value = 63600
$ws.Range("N71").Value = -198912
$ws.Range("L71").Value = 190800
$ws.Range("J71").Value = 63600
$ws.Range("J110").Value = 1893.3334
$ws.Range("N110").Value = -5983.3334
$ws.Range("M110").Value = 924.0714
$ws.Range("K110").Value = 1120.9286
$ws.Range("I110").Value = 1120.9286
$ws.Range("H110").Value = 1257.2354
$ws.Range("L110").Value = 1893.3334
$ws.Range("H132").Value = 823099.5
$ws.Range("L132").Value = 9342
$ws.Range("K132").Value = 3048111.9
$ws.Range("M132").Value = -3045581.9
$ws.Range("I132").Value = 1016037.3
$ws.Range("N132").Value = -14402
$ws.Range("J132").Value = 3114

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 765.0625
$ws.Range("L80").Value = 1146.8889
$ws.Range("I80").Value = 274.14285
$ws.Range("K80").Value = 274.14285
$ws.Range("J80").Value = 1146.8889
$ws.Range("N80").Value = -3142.8889
$ws.Range("M80").Value = 723.85715
$ws.Range("I83").Value = 274.14285
$ws.Range("J83").Value = 1146.8889
$ws.Range("H83").Value = 765.0625
$ws.Range("N83").Value = -15718.4445
$ws.Range("M83").Value = 3621.28575
$ws.Range("K83").Value = 1370.71425
$ws.Range("L83").Value = 5734.4445

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("M31").Value = -929.1538
$ws.Range("K31").Value = 1224.1538
$ws.Range("N31").Value = -4090
$ws.Range("H31").Value = 1527.6
$ws.Range("L31").Value = 3500
$ws.Range("I31").Value = 1224.1538
$ws.Range("J31").Value = 3500
$ws.Range("M34").Value = -1022.1538
$ws.Range("K34").Value = 1224.1538
$ws.Range("H34").Value = 1527.6
$ws.Range("N34").Value = -3904
$ws.Range("L34").Value = 3500
$ws.Range("J34").Value = 3500
$ws.Range("I34").Value = 1224.1538
$ws.Range("H132").Value = 2414.4194
$ws.Range("L132").Value = 8096.1819
$ws.Range("K132").Value = 6774.150000000001
$ws.Range("M132").Value = -4244.150000000001
$ws.Range("I132").Value = 2258.05
$ws.Range("N132").Value = -13156.1819
$ws.Range("J132").Value = 2698.7273

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 15625652
$ws.Range("J5").Value = 55556444
$ws.Range("L5").Value = 166669332
$ws.Range("K5").Value = 1678.5651
$ws.Range("M5").Value = -1566.5651
$ws.Range("I5").Value = 559.5217
$ws.Range("N5").Value = -166669556
$ws.Range("K70").Value = 17778366
$ws.Range("M70").Value = -17778051
$ws.Range("I70").Value = 5926122
$ws.Range("J70").Value = 3980
$ws.Range("L70").Value = 11940
$ws.Range("H70").Value = 4445586.5
$ws.Range("N70").Value = -12570
$ws.Range("M73").Value = -17777274
$ws.Range("L73").Value = 11940
$ws.Range("I73").Value = 5926122
$ws.Range("N73").Value = -14124
$ws.Range("K73").Value = 17778366
$ws.Range("H73").Value = 4445586.5
$ws.Range("J73").Value = 3980
$ws.Range("I135").Value = 559.5217
$ws.Range("K135").Value = 5035.6953
$ws.Range("H135").Value = 15625652
$ws.Range("M135").Value = -2500.6953
$ws.Range("L135").Value = 500007996
$ws.Range("J135").Value = 55556444
$ws.Range("N135").Value = -500013066

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1077.2354
$ws.Range("K102").Value = 894.3077
$ws.Range("M102").Value = 727.6923
$ws.Range("I102").Value = 894.3077
$ws.Range("J126").Value = 2339.8
$ws.Range("I126").Value = 1637.3334
$ws.Range("K126").Value = 4912.0002
$ws.Range("M126").Value = -2442.0002
$ws.Range("N126").Value = -11959.4
$ws.Range("H126").Value = 2076.375
$ws.Range("L126").Value = 7019.400000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("L7").Value = 2139.5
$ws.Range("N7").Value = -2363.5
$ws.Range("I7").Value = 1547.2858
$ws.Range("K7").Value = 1547.2858
$ws.Range("M7").Value = -1435.2858
$ws.Range("H7").Value = 1703.1316
$ws.Range("J7").Value = 2139.5
$ws.Range("M122").Value = -13617.5005
$ws.Range("K122").Value = 16067.5005
$ws.Range("H122").Value = 4866.622
$ws.Range("I122").Value = 5355.8335
$ws.Range("J126").Value = 2139.5
$ws.Range("I126").Value = 1547.2858
$ws.Range("K126").Value = 4641.857400000001
$ws.Range("M126").Value = -2171.857400000001
$ws.Range("N126").Value = -11358.5
$ws.Range("H126").Value = 1703.1316
$ws.Range("L126").Value = 6418.5
